# Refresh the cryptocurrency Price (D) / Volume(1h) (E) columns with the
# latest scraped figures. Both columns hold plain text in the workbook
# (e.g. "29.423.53", "  -0.34%  "), so for any new Price value that would
# otherwise be auto-recognized as a number by Excel, the cell is briefly
# switched to Text format, assigned, then had its explicit formatting
# cleared again so it keeps the original (unstyled) look while staying
# text-typed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.423.53'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '1.849.11'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.94'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6345'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.787.96'
$ws.Range("E8").Value = '  +102.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07561'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("E10").Value = '  -1.00%  '
$ws.Range("D11").Value = '4.013.32'
$ws.Range("E11").Value = '  +88.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.65'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07717'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.994'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6864'
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.03'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009930'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.212'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.40%  '
$ws.Range("D19").Value = '29.447.45'
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '232.07'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.50'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.603'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.28%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.83'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1387'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.423'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.69'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("D29").Value = '3.977.79'
$ws.Range("E29").Value = '  +97.73%  '
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05816'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.260'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("E34").Value = '  -1.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.861'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.158'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7170'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.589'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.34%  '
$ws.Range("D39").Value = '1.252.19'
$ws.Range("E39").Value = '  +4.05%  '
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01808'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9014'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.112'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.72'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '67.09'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.211'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.154'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4019'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.687'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("E51").Value = '  -0.04%  '
